# Update matchDay dates (column A) based on round 1 results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @{
    2  = 45732
    3  = 45731
    4  = 45729
    5  = 45731
    6  = 45730
    7  = 45731
    8  = 45732
    9  = 45731
    10 = 45732
    11 = 45730
    12 = 45732
    13 = 45731
    14 = 45731
    15 = 45729
    16 = 45732
    17 = 45731
    18 = 45732
    19 = 45731
}

foreach ($row in $newDates.Keys) {
    $ws.Cells.Item($row, 1).Value2 = $newDates[$row]
}
